$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Full results" ---
$ws1 = $wb.Worksheets.Item("Full results")

# Row 2 (wealth / NULL MODEL)
$ws1.Range("C2").Value = 0.891438339859624
$ws1.Range("D2").Value = 0.108676415936547
$ws1.Range("E2").Value = 1.00011475579617
$ws1.Range("J2").Value = 0.108663946118895
$ws1.Range("K2").Value = 0.0645677439092995
$ws1.Range("L2").Value = -0.0109867315004375
$ws1.Range("M2").Value = -0.0406373313050686
$ws1.Range("N2").Value = 0.053581012408862

# Row 3 (wealth / CONDITIONAL MODEL)
$ws1.Range("F3").Value = 0.921092342242442
$ws1.Range("G3").Value = 0.0645751534321587

# Row 4 (wealth / COMPLETE MODEL)
$ws1.Range("H4").Value = 0.932080334534
$ws1.Range("I4").Value = 0.0156529172020237
$ws1.Range("O4").Value = 0.0680266148138263

# --- Sheet 2: "For plotting" ---
$ws2 = $wb.Worksheets.Item("For plotting")

# Row 2 (Sibcorr / wealth)
$ws2.Range("C2").Value = 0.108663946118895
$ws2.Range("D2").Value = 0.0492859953086397
$ws2.Range("E2").Value = 0.16804189692915
$ws2.Range("F2").Value = 948

# Row 3 (IOLIB / wealth)
$ws2.Range("C3").Value = 0.053581012408862
$ws2.Range("D3").Value = -0.00557082496065574
$ws2.Range("E3").Value = 0.11273284977838
$ws2.Range("F3").Value = 948

# Row 4 (IORAD / wealth)
$ws2.Range("C4").Value = 0.0680266148138263
$ws2.Range("D4").Value = 0.00910844429976228
$ws2.Range("E4").Value = 0.12694478532789
$ws2.Range("F4").Value = 948
